$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "98.098.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.17%  "

$ws.Range("D3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.378.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.05%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.24%  "

$ws.Range("D6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "657.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.61%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.426"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.27%  "

$ws.Range("E9").Value = "  -0.07%  "

$ws.Range("E10").Value = "  -3.43%  "

$ws.Range("D11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.375.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.06%  "

$ws.Range("E12").Value = "  -3.24%  "

$ws.Range("D13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "43.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.99%  "

$ws.Range("D14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "97.841.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.21%  "

$ws.Range("D15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.11%  "

$ws.Range("D16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000257"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.65%  "

$ws.Range("D17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.017.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.81%  "

$ws.Range("D18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.05%  "

$ws.Range("D19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.396.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.38%  "

$ws.Range("D20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.14%  "

$ws.Range("D21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.518"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -11.85%  "

$ws.Range("D22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.10%  "

$ws.Range("D23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "508.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.57%  "

$ws.Range("E24").Value = "  -1.02%  "

$ws.Range("E25").Value = "  -2.89%  "

$ws.Range("D26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.28%  "

$ws.Range("D27").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "93.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.83%  "

$ws.Range("D28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.67%  "

$ws.Range("D29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.556.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.30%  "

$ws.Range("D30").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.35%  "

$ws.Range("D31").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.143"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.01%  "

$ws.Range("E32").Value = "  -0.13%  "

$ws.Range("E33").Value = "  -3.91%  "

$ws.Range("D34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.74%  "

$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("E36").Value = "  -2.61%  "

$ws.Range("D37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "28.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.40%  "

$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.16%  "

$ws.Range("D40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "523.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.47%  "

$ws.Range("E41").Value = "  -1.36%  "

$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("E43").Value = "  -1.17%  "

$ws.Range("D44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.55%  "

$ws.Range("D45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.848"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.26%  "

$ws.Range("D46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0427"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.03%  "

$ws.Range("E47").Value = "  -3.02%  "

$ws.Range("E48").Value = "  +6.96%  "

$ws.Range("D49").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.53%  "

$ws.Range("E50").Value = "  -4.51%  "

$ws.Range("D51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.30%  "
